# New order #24 came in after the 2026-01-20 15:03 snapshot was taken.
# It gets prepended to the "All Orders" log (newest first) and the
# "Daily Summary" rollup for 2026-01-20 is refreshed to include it.

$wb = $excel.ActiveWorkbook

# --- "All Orders": insert the new order as the new top data row ---
$ws1 = $wb.Worksheets.Item("All Orders")
$ws1.Rows("2:2").Insert()

$ws1.Cells.Item(2, 1).Value = 24
$ws1.Cells.Item(2, 2).Value = "2026-01-20 15:03"
$ws1.Cells.Item(2, 3).Value = "Udita Roy"
$ws1.Cells.Item(2, 4).Value = "A-1603"

# Phone numbers are kept as text (leading-zero safe) like the rest of column E.
$ws1.Cells.Item(2, 5).NumberFormat = "@"
$ws1.Cells.Item(2, 5).Value = "7061856166"

$ws1.Cells.Item(2, 6).Value = "Appe Chutney x1"
$ws1.Cells.Item(2, 7).Value = 60
$ws1.Cells.Item(2, 8).Value = "NEW"
$ws1.Cells.Item(2, 9).Value = "PENDING"

# Collection Date is kept as plain text (YYYY-MM-DD) like the rest of column J,
# not auto-converted to a date serial.
$ws1.Cells.Item(2, 10).NumberFormat = "@"
$ws1.Cells.Item(2, 10).Value = "2026-01-21"

$ws1.Cells.Item(2, 11).Value = "09:30"
$ws1.Cells.Item(2, 12).Value = "Less spicy. Flavourful"

# --- "Daily Summary": refresh the 2026-01-20 rollup for the new order ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(2, 2).Value = 5    # Total Orders: 4 -> 5
$ws2.Cells.Item(2, 5).Value = 320  # Revenue: 260 -> 320
$ws2.Cells.Item(2, 7).Value = 270  # Pending: 210 -> 270
